$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# Insert two new columns before column L (existing "Bet Size" header), so the
# existing L:S content shifts right to N:U, matching the target layout.
$insertRange = $ws.Range("L1:M1").EntireColumn
$insertRange.Insert()

# New column K: "Expected Payout" header + per-row expected payout formula.
$ws.Range("K1").Value = "Expected Payout"
$ws.Range("K2").Formula = '=$O$1*SUM(C2*1.5,D2,E2*2,-G2*2,-H2,-I2*2,-J2)'

# New column L: "Validator" header + row-to-row reconciliation check.
$ws.Range("L1").Value = "Validator"
$ws.Range("L2").Formula = '=IF(K2+B2=B3,"","XXXX")'

# Size the two new columns to fit their (now longer) header text, matching
# the "Expected Payout" / "Validator" column widths from the authored sheet.
$ws.Columns.Item(11).ColumnWidth = 14.17
$ws.Columns.Item(12).ColumnWidth = 7.6
